# repull data, push all data, mean calculation
# Update the dSF (column F) values to reflect the re-pulled/recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    5  = -1
    7  = -2
    10 = 1
    14 = -3
    15 = 0
    16 = -2
    18 = -3
    20 = -4
    21 = 1
    23 = 3
    24 = 4
    25 = -1
    33 = -1
    38 = 0
    41 = 3
    43 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
